# Update the division-problem table cells to the new values described
# by the commit diff. The worksheet/document contains a single 5-column
# table; the data rows (1-indexed) are 1, 5, 9, 13 and 17 - the rows in
# between are blank spacer rows. We address each cell positionally so
# duplicate text values (several cells share identical old/new text)
# are never ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row => ordered list of new cell values (columns 1..5)
$updates = @{
    1  = @("60÷3=20, 0", "50÷7=7, 1", "93÷8=11, 5", "56÷8=7, 0", "26÷8=3, 2")
    5  = @("69÷7=9, 6", "62÷3=20, 2", "52÷3=17, 1", "75÷9=8, 3", "38÷6=6, 2")
    9  = @("33÷5=6, 3", "81÷6=13, 3", "55÷3=18, 1", "93÷9=10, 3", "60÷5=12, 0")
    13 = @("96÷7=13, 5", "46÷5=9, 1", "25÷3=8, 1", "56÷2=28, 0", "60÷5=12, 0")
    17 = @("30÷6=5, 0", "30÷2=15, 0", "88÷4=22, 0", "82÷9=9, 1", "16÷8=2, 0")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
